$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(19, 37, 56, 2023.455715725201, -1757.143969265623, 0.2035205594775814, -0.174833651311214, -0.04646347265072986, 38445.65859877883, -65014.32686282805, 0.3392857142857143, 1.151559434581152, 0.5913413312714022, -26568.66826404922)
    3 = @(76, 128, 204, 2659.122434410413, -1266.878094258235, 0.2653458717649267, -0.1266188768220287, 0.01940720598487632, 202093.3050151914, -162160.3960650541, 0.3725490196078431, 2.098956834491124, 1.246255620479105, 39932.90895013732)
    4 = @(459, 252, 711, 3515.771750332336, -1461.611037167878, 0.3506835221172528, -0.1450228047352465, 0.1749901404480124, 1613739.233402542, -368325.9813663051, 0.6455696202531646, 2.405408594303411, 4.381279939624071, 1245413.252036237)
    5 = @(169, 226, 395, 3985.434833530434, -1478.200596174383, 0.3991966649029029, -0.145709615820076, 0.08742750175507198, 673538.486866643, -334073.3347354103, 0.4278481012658228, 2.69613937637749, 2.016139622158389, 339465.1521312327)
    6 = @(151, 162, 313, 3396.930002765135, -1465.858696723967, 0.3354916368576785, -0.1448727135832827, 0.08686855452082315, 512936.4304175353, -237469.1088692828, 0.4824281150159744, 2.317365248340035, 2.160013287032994, 275467.3215482525)
    7 = @(133, 155, 288, 3336.579013530943, -1294.467023875952, 0.343679563965126, -0.1302561843829433, 0.08860997718057478, 443765.0087996155, -200642.3887007725, 0.4618055555555556, 2.577569727145623, 2.21172112071205, 243122.6200988429)
    8 = @(97, 166, 263, 2893.18213236662, -1484.115858605109, 0.2883639850999758, -0.1487136588552809, 0.01248988283163889, 280638.6668395621, -246363.232528448, 0.3688212927756654, 1.94943145145411, 1.139125607174992, 34275.43431111405)
    9 = @(19, 28, 47, 1225.213359007508, -1324.689633966131, 0.1166077775345663, -0.1349916419456071, -0.03328123832596255, 23279.05382114265, -37091.30975105168, 0.4042553191489361, 0.9249059761562483, 0.6276147695345969, -13812.25592990903)
    10 = @(317, 370, 687, 3487.484046433025, -1423.738000469894, 0.3468196139453806, -0.1417020697746878, 0.08371477700735255, 1105532.44271927, -526783.060173861, 0.4614264919941776, 2.449526559860036, 2.098648431015221, 578749.382545409)
    11 = @(240, 281, 521, 2703.860106109207, -1409.904036393012, 0.2700878246327166, -0.1394988405077706, 0.04917831809821207, 648926.4254662098, -396183.0342264363, 0.4606525911708254, 1.91776180244618, 1.637946023441578, 252743.3912397735)
    12 = @(95, 179, 274, 3692.408231267399, -1478.119580517253, 0.3796808991826974, -0.1474264096549064, 0.03532977406616068, 350778.7819704029, -264583.4049125884, 0.3467153284671533, 2.498044326004583, 1.325777714918633, 86195.37705781456)
    13 = @(252, 292, 544, 4061.933210393735, -1532.938010490335, 0.4077591433576939, -0.1532105081895778, 0.1066504333359967, 1023607.169019221, -447617.899063178, 0.4632352941176471, 2.649770038055525, 2.286787841061616, 575989.2699560431)
    14 = @(439, 333, 772, 2949.995968620971, -1371.826974241163, 0.2954755674806395, -0.1373189026705082, 0.1087909061330591, 1295048.230224607, -456818.3824223069, 0.5686528497409327, 2.150414027434317, 2.834930204335334, 838229.8478023004)
    15 = @(143, 245, 388, 2785.829363638858, -1410.398998185829, 0.2771690347713198, -0.1403765942487842, 0.01351264531274894, 398373.5990003567, -345547.7545555281, 0.3685567010309279, 1.975206567235386, 1.152875669855756, 52825.84444482869)
    16 = @(101, 169, 270, 2542.021558192822, -1650.468174679579, 0.255111345037014, -0.1645412831531281, -0.007560114830149, 256744.177377475, -278929.1215208489, 0.3740740740740741, 1.540182111470479, 0.9204638654350202, -22184.94414337384)
    17 = @(377, 389, 766, 2499.623027194526, -1466.390514183418, 0.2502988253769945, -0.1454475982551954, 0.04932577212252722, 942357.8812523355, -570425.91001735, 0.4921671018276763, 1.704609381346468, 1.652025030250945, 371931.9712349855)
    18 = @(315, 404, 719, 2714.222590407798, -1483.622255727249, 0.2707503044128109, -0.1478941101394559, 0.0355175596574344, 854980.1159784566, -599383.3913138085, 0.4381084840055633, 1.829456642302341, 1.426432778032766, 255596.7246646481)
    19 = @(211, 425, 636, 3250.501439966892, -1494.710566602336, 0.3212071353777778, -0.1493090841011612, 0.00678985034861256, 685855.8038330143, -635251.9908059925, 0.3317610062893082, 2.174669472870383, 1.079659432413296, 50603.81302702182)
    20 = @(112, 314, 426, 3142.936185147014, -1636.154498776297, 0.3134494740534706, -0.1629852793269602, -0.0377254380626216, 352008.8527364655, -513752.5126157572, 0.2629107981220657, 1.92092873105667, 0.685172031459704, -161743.6598792917)
    21 = @(198, 514, 712, 3486.361317093676, -1718.42979870873, 0.3462711062050057, -0.1688371728405392, -0.02559076939809819, 690299.5407845477, -883272.9165362868, 0.2780898876404495, 2.028806367134353, 0.7815246316976692, -192973.3757517391)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $arr = New-Object "object[,]" 1,14
    for ($i = 0; $i -lt 14; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range("C$row`:P$row").Value = $arr
}
